# Auto-generated edit script: updates Horarios Linea 141 schedules
# across sheets LP1912, LP1912-215 and 6203-6173 (commit: "Horarios actualizados Linea 141 - 878")
$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 ----
$ws = $wb.Worksheets.Item("LP1912")

$ws.Range("A1").Value = "LÍNEA 141 - LP1912 - 24/01/2026"
$ws.Range("A2").Value = "Última actualización: 04:44:46"
$ws.Range("A3").Value = "Total filas: 33"

# Data rows 5 (column labels) through 38
$rows = @(
    @("Hora_Scrap","Hora_Llegada","Linea","Minutos","Parada"),
    @("03:52:04","04:01","81_EL PELIGRO",9,"LP1912"),
    @("04:32:18","04:33","15_ABASTO",1,"LP1912"),
    @("04:44:46","04:46","15_ABASTO",2,"LP1912"),
    @("04:44:46","04:46","215_EL PELIGRO",2,"LP1912"),
    @("03:52:04","04:46","215A_EL PATO",54,"LP1912"),
    @("04:32:18","04:47","215_EL PELIGRO",15,"LP1912"),
    @("04:44:46","04:53","11_ETCHEVERRY",9,"LP1912"),
    @("04:13:31","05:11","17_ROMERO",58,"LP1912"),
    @("03:52:04","05:16","17_ROMERO",84,"LP1912"),
    @("04:44:46","05:22","23_HERNANDEZ",38,"LP1912"),
    @("04:44:46","05:31","81_EL PELIGRO",47,"LP1912"),
    @("04:32:18","05:32","81_EL PELIGRO",60,"LP1912"),
    @("03:52:04","05:35","215B_EL PATO",103,"LP1912"),
    @("04:44:46","05:44","14_ABASTO",60,"LP1912"),
    @("03:52:04","05:46","15_ABASTO",114,"LP1912"),
    @("04:32:18","05:47","14_ABASTO",75,"LP1912"),
    @("04:13:31","05:50","14_ABASTO",97,"LP1912"),
    @("04:44:46","05:51","17_ROMERO",67,"LP1912"),
    @("04:32:18","05:52","17_ROMERO",80,"LP1912"),
    @("04:44:46","06:00","16_SANTA ANA",76,"LP1912"),
    @("04:32:18","06:01","16_SANTA ANA",89,"LP1912"),
    @("04:44:46","06:03","10_OLMOS",79,"LP1912"),
    @("04:32:18","06:04","10_OLMOS",92,"LP1912"),
    @("04:44:46","06:10","215A_EL PATO",86,"LP1912"),
    @("04:32:18","06:11","215A_EL PATO",99,"LP1912"),
    @("04:32:18","06:15","17_ROMERO",103,"LP1912"),
    @("04:44:46","06:24","11_ETCHEVERRY",100,"LP1912"),
    @("04:44:46","06:27","23_HERNANDEZ",103,"LP1912"),
    @("04:44:46","06:28","17_ROMERO",104,"LP1912"),
    @("04:44:46","06:30","16_SANTA ANA",106,"LP1912"),
    @("04:32:18","06:31","16_SANTA ANA",119,"LP1912"),
    @("04:44:46","06:31","17X38_ROMERO",107,"LP1912"),
    @("04:44:46","06:39","225_C ROCA-H SUR",115,"LP1912")
)

$startRow = 5
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    if ($row[3] -ne $null) { $ws.Cells.Item($r, 4).Value = $row[3] }
    $ws.Cells.Item($r, 5).Value = $row[4]
}

# ---- Sheet: LP1912-215 ----
$ws = $wb.Worksheets.Item("LP1912-215")

$ws.Range("A1").Value = "LÍNEA 141 - LP1912-215 - 24/01/2026"
$ws.Range("A2").Value = "Última actualización: 04:44:46"
$ws.Range("A3").Value = "Total filas: 6"

# Data rows 5 (column labels) through 11
$rows = @(
    @("Hora_Scrap","Hora_Llegada","Linea","Minutos","Parada"),
    @("04:44:46","04:46","215_EL PELIGRO",2,"LP1912"),
    @("03:52:04","04:46","215A_EL PATO",54,"LP1912"),
    @("04:32:18","04:47","215_EL PELIGRO",15,"LP1912"),
    @("03:52:04","05:35","215B_EL PATO",103,"LP1912"),
    @("04:44:46","06:10","215A_EL PATO",86,"LP1912"),
    @("04:32:18","06:11","215A_EL PATO",99,"LP1912")
)

$startRow = 5
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    if ($row[3] -ne $null) { $ws.Cells.Item($r, 4).Value = $row[3] }
    $ws.Cells.Item($r, 5).Value = $row[4]
}

# ---- Sheet: 6203-6173 ----
$ws = $wb.Worksheets.Item("6203-6173")

$ws.Range("A1").Value = "LÍNEA 141 - 6203-6173 - 24/01/2026"
$ws.Range("A2").Value = "Última actualización: 04:44:46"
$ws.Range("A3").Value = "Total filas: 1"

# Data rows 5 (column labels) through 6
$rows = @(
    @("Hora_Scrap","Hora_Llegada","Linea","Minutos","Parada"),
    @("03:52:04","05:44","215A_LA PLATA",112,"L6173")
)

$startRow = 5
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    if ($row[3] -ne $null) { $ws.Cells.Item($r, 4).Value = $row[3] }
    $ws.Cells.Item($r, 5).Value = $row[4]
}
